$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'256.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'4.70%"
$ws.Range("E2").Style = "Normal"

$ws.Range("E3").Value = "'-2.82%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'5.209"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.43%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.05911"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.50%"
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'0.41%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.8658"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.86%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'1.007"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'14.24%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.1415"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'2.40%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.03579"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'7.97%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07218"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.92%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.03176"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.02%"
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'0.10%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.001548"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'1.44%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.0006039"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-93.99%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.005869"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-2.85%"
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'3.478"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.53%"
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'1.84%"
$ws.Range("E18").Style = "Normal"

$ws.Range("E20").Value = "'-0.50%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'0.1308"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.07%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'3.550"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.72%"
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'2.54%"
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'1.48%"
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.001220"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.12%"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.004518"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'8.72%"
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'0.05%"
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'33.80%"
$ws.Range("E28").Style = "Normal"

$ws.Range("D40").Value = "'0.03818"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.84%"
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = "'KickToken"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.005600"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'8.77%"
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = "'BKEXToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.1101"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'3.37%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.001900"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-15.51%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.01067"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'12.68%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.00005427"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'3.04%"
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'0.05%"
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'22.49%"
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'-4.12%"
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'0.05%"
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'0.05%"
$ws.Range("E50").Style = "Normal"
